# YOLOX Custom Training Guide — apply commit changes
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "( If CUDA related ... )" paragraph: extend text and make whole
#    paragraph red (FF0000), matching the new "for different CUDA version"
#    wording and highlighting it as an important note.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "website )", $true, $false, $false, $false, $false,
    $true, 1, $false, "website for different CUDA version )", 2) | Out-Null

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*( If*CUDA related*") {
        $p.Range.Font.Color = 255
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Point readers at the exact default config folder inside the repo.
# ---------------------------------------------------------------------------
$old3 = "Create a new experiment configuration file based on ``yolox_"
$new3 = "Create a new experiment configuration file based on ``YOLOX/exp/defaults/yolox_"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) New paragraph describing how to copy/rename the experiment file,
#    inserted right after the "...hyperparameters." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*and customize paths, number of classes, and hyperparameters.*") {
        $p.Range.InsertParagraphAfter()
        $newP = $p.Next()
        $newP.Range.Text = "Copy the “yolox.py” and in the same folder paste it and rename it to anything like “my_exp.py”. Now go to the main GitHub page and see what the best configuration setup for training is. Copy it and make necessary changes according to your needs."
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Training command gains a --cache flag.
# ---------------------------------------------------------------------------
$oldTrain = "python tools/train.py -f exps/default/exp_my_dataset.py -d 1 -b 8 --fp16 -c yolox_s.pth"
$newTrain = "python tools/train.py -f exps/default/exp_my_dataset.py -d 1 -b 8 --fp16 -c yolox_s.pth --cache"
$d.Content.Find.Execute($oldTrain, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newTrain, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) After the (updated) training command, insert:
#      - a red "(Note:- ... yolox_s.pth ... )" callout
#      - guidance + command for resuming interrupted training
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*--fp16 -c yolox_s.pth --cache*") {

        # Insert all four paragraphs first (while still plain/black), then
        # colour only the "Note" paragraph — this avoids Word's new-paragraph
        # formatting inheritance from leaking the red colour into the rest.
        $p.Range.InsertParagraphAfter()
        $noteP = $p.Next()
        $noteP.Range.Text = "( Note:- there is a file named “yolox_s.pth” at last, of the above command. It is a pre trained weight you need to download from the yolox official github page.)"

        $noteP.Range.InsertParagraphAfter()
        $resumeIntroP = $noteP.Next()
        $resumeIntroP.Range.Text = "If in future your training gets terminated in middle then you could resume the training using the following command:"

        $resumeIntroP.Range.InsertParagraphAfter()
        $bashP = $resumeIntroP.Next()
        $bashP.Range.Text = "bash"

        $bashP.Range.InsertParagraphAfter()
        $resumeCmdP = $bashP.Next()
        $resumeCmdP.Range.Text = "python tools/train.py -f exps/default/base_exp.py -d 2 -b 8 --resume --cache"

        $noteP.Range.Font.Color = 255

        break
    }
}

# ---------------------------------------------------------------------------
# 6) Drop the stray lastRenderedPageBreak markers that moved elsewhere in
#    the document (rewriting the run text clears them cleanly).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*TensorRT:**") {
        $p.Range.Text = "**TensorRT:**"
        break
    }
}

$videoSectionSeen = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Single Video:**") {
        $videoSectionSeen = $true
        continue
    }
    if ($videoSectionSeen -and ($p.Range.Text -like "*bash*")) {
        $p.Range.Text = "bash"
        break
    }
}
